$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.586.64"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "1.924.67"
$ws.Range("E3").Value = "  +3.55%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'247.07"
$ws.Range("E5").Value = "  +4.61%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "'0.4737"
$ws.Range("E7").Value = "  +1.76%  "
$ws.Range("D8").Value = "'0.2912"
$ws.Range("E8").Value = "  +3.50%  "
$ws.Range("D9").Value = "'0.06784"
$ws.Range("E9").Value = "  +6.26%  "
$ws.Range("D10").Value = "'104.89"
$ws.Range("E10").Value = "  +9.77%  "
$ws.Range("D11").Value = "'18.44"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").Value = "1.912.86"
$ws.Range("E12").Value = "  +2.86%  "
$ws.Range("D13").Value = "'0.07726"
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("D14").Value = "'5.327"
$ws.Range("E14").Value = "  +7.18%  "
$ws.Range("D15").Value = "'0.6734"
$ws.Range("E15").Value = "  +5.27%  "
$ws.Range("D16").Value = "'287.45"
$ws.Range("E16").Value = "  -2.51%  "
$ws.Range("D17").Value = "30.612.91"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("D18").Value = "'0.000007632"
$ws.Range("E18").Value = "  +3.04%  "
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "2.164.50"
$ws.Range("E21").Value = "  +2.61%  "
$ws.Range("D22").Value = "'5.451"
$ws.Range("E22").Value = "  +9.02%  "
$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'6.313"
$ws.Range("E24").Value = "  +4.61%  "
$ws.Range("D25").Value = "'9.406"
$ws.Range("E25").Value = "  +4.04%  "
$ws.Range("D26").Value = "'168.24"
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("D27").Value = "'20.82"
$ws.Range("E27").Value = "  +7.76%  "
$ws.Range("D28").Value = "'2.150"
$ws.Range("E28").Value = "  +11.12%  "
$ws.Range("D29").Value = "'0.1085"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").Value = "'1.362"
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("D31").Value = "'4.193"
$ws.Range("E31").Value = "  +3.88%  "
$ws.Range("D32").Value = "'4.203"
$ws.Range("E32").Value = "  +10.52%  "
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("D34").Value = "'0.7427"
$ws.Range("E34").Value = "  +2.50%  "
$ws.Range("D35").Value = "'1.164"
$ws.Range("E35").Value = "  +4.15%  "
$ws.Range("E36").Value = "  +7.52%  "
$ws.Range("D37").Value = "'2.746"
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("D38").Value = "'2.694"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").Value = "'2.064"
$ws.Range("E39").Value = "  +4.77%  "
$ws.Range("D40").Value = "'111.81"
$ws.Range("E40").Value = "  +5.82%  "
$ws.Range("D41").Value = "'0.8839"
$ws.Range("E41").Value = "  +2.20%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.4386"
$ws.Range("E42").Value = "  +8.04%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.959"
$ws.Range("E43").Value = "  +6.61%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").Value = "'67.42"
$ws.Range("E45").Value = "  +3.03%  "
$ws.Range("D46").Value = "'7.278"
$ws.Range("E46").Value = "  +2.77%  "
$ws.Range("D47").Value = "'9.321"
$ws.Range("E47").Value = "  +3.45%  "
$ws.Range("D48").Value = "'48.13"
$ws.Range("E48").Value = "  +16.98%  "
$ws.Range("D49").Value = "'0.1235"
$ws.Range("E49").Value = "  +3.90%  "
$ws.Range("D50").Value = "'35.24"
$ws.Range("E50").Value = "  +4.16%  "
$ws.Range("D51").Value = "'0.4058"
$ws.Range("E51").Value = "  +8.61%  "

# Reset style on cells whose new text looks numeric, since setting
# them via an apostrophe-prefix (to keep them as text, matching the
# original inlineStr cell type) adds Excel quote-prefix styling;
# resetting to the Normal style keeps formatting identical to before.
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
